$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the full target range is formatted as Text so numeric-looking values
# (scores, balls, strike rates) are preserved as strings, matching the source data.
$ws.Range("A1:K15").NumberFormat = "@"

# Row 1
$ws.Cells.Item(1, 1).Value = "venue"
$ws.Cells.Item(1, 2).Value = "date"
$ws.Cells.Item(1, 3).Value = "result"
$ws.Cells.Item(1, 4).Value = "ownTeam"
$ws.Cells.Item(1, 5).Value = "oppTeam"
$ws.Cells.Item(1, 6).Value = "batsman"
$ws.Cells.Item(1, 7).Value = "totalRuns"
$ws.Cells.Item(1, 8).Value = "totalBalls"
$ws.Cells.Item(1, 9).Value = "total4s"
$ws.Cells.Item(1, 10).Value = "total6s"
$ws.Cells.Item(1, 11).Value = "sr"

# Row 2
$ws.Cells.Item(2, 1).Value = " Abu Dhabi"
$ws.Cells.Item(2, 2).Value = " October 07 2020"
$ws.Cells.Item(2, 3).Value = "KKR won by 10 runs"
$ws.Cells.Item(2, 4).Value = "Kolkata Knight Riders"
$ws.Cells.Item(2, 5).Value = "Chennai Super Kings"
$ws.Cells.Item(2, 6).Value = "Shubman Gill "
$ws.Cells.Item(2, 7).Value = "11"
$ws.Cells.Item(2, 8).Value = "12"
$ws.Cells.Item(2, 9).Value = "1"
$ws.Cells.Item(2, 10).Value = "0"
$ws.Cells.Item(2, 11).Value = "91.66"

# Row 3
$ws.Cells.Item(3, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(3, 2).Value = " September 30 2020"
$ws.Cells.Item(3, 3).Value = "KKR won by 37 runs"
$ws.Cells.Item(3, 4).Value = "Kolkata Knight Riders"
$ws.Cells.Item(3, 5).Value = "Rajasthan Royals"
$ws.Cells.Item(3, 6).Value = "Shubman Gill "
$ws.Cells.Item(3, 7).Value = "47"
$ws.Cells.Item(3, 8).Value = "34"
$ws.Cells.Item(3, 9).Value = "5"
$ws.Cells.Item(3, 10).Value = "1"
$ws.Cells.Item(3, 11).Value = "138.23"

# Row 4
$ws.Cells.Item(4, 1).Value = " Abu Dhabi"
$ws.Cells.Item(4, 2).Value = " October 16 2020"
$ws.Cells.Item(4, 3).Value = "Mumbai won by 8 wickets (with 19 balls remaining)"
$ws.Cells.Item(4, 4).Value = "Kolkata Knight Riders"
$ws.Cells.Item(4, 5).Value = "Mumbai Indians"
$ws.Cells.Item(4, 6).Value = "Shubman Gill "
$ws.Cells.Item(4, 7).Value = "21"
$ws.Cells.Item(4, 8).Value = "23"
$ws.Cells.Item(4, 9).Value = "2"
$ws.Cells.Item(4, 10).Value = "0"
$ws.Cells.Item(4, 11).Value = "91.30"

# Row 5
$ws.Cells.Item(5, 1).Value = " Abu Dhabi"
$ws.Cells.Item(5, 2).Value = " October 18 2020"
$ws.Cells.Item(5, 3).Value = "Match tied (KKR won the one-over eliminator)"
$ws.Cells.Item(5, 4).Value = "Kolkata Knight Riders"
$ws.Cells.Item(5, 5).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(5, 6).Value = "Shubman Gill "
$ws.Cells.Item(5, 7).Value = "36"
$ws.Cells.Item(5, 8).Value = "37"
$ws.Cells.Item(5, 9).Value = "5"
$ws.Cells.Item(5, 10).Value = "0"
$ws.Cells.Item(5, 11).Value = "97.29"

# Row 6
$ws.Cells.Item(6, 1).Value = " Abu Dhabi"
$ws.Cells.Item(6, 2).Value = " October 10 2020"
$ws.Cells.Item(6, 3).Value = "KKR won by 2 runs"
$ws.Cells.Item(6, 4).Value = "Kolkata Knight Riders"
$ws.Cells.Item(6, 5).Value = "Kings XI Punjab"
$ws.Cells.Item(6, 6).Value = "Shubman Gill "
$ws.Cells.Item(6, 7).Value = "57"
$ws.Cells.Item(6, 8).Value = "47"
$ws.Cells.Item(6, 9).Value = "5"
$ws.Cells.Item(6, 10).Value = "0"
$ws.Cells.Item(6, 11).Value = "121.27"

# Row 7
$ws.Cells.Item(7, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(7, 2).Value = " November 01 2020"
$ws.Cells.Item(7, 3).Value = "KKR won by 60 runs"
$ws.Cells.Item(7, 4).Value = "Kolkata Knight Riders"
$ws.Cells.Item(7, 5).Value = "Rajasthan Royals"
$ws.Cells.Item(7, 6).Value = "Shubman Gill "
$ws.Cells.Item(7, 7).Value = "36"
$ws.Cells.Item(7, 8).Value = "24"
$ws.Cells.Item(7, 9).Value = "6"
$ws.Cells.Item(7, 10).Value = "0"
$ws.Cells.Item(7, 11).Value = "150.00"

# Row 8
$ws.Cells.Item(8, 1).Value = " Abu Dhabi"
$ws.Cells.Item(8, 2).Value = " September 26 2020"
$ws.Cells.Item(8, 3).Value = "KKR won by 7 wickets (with 12 balls remaining)"
$ws.Cells.Item(8, 4).Value = "Kolkata Knight Riders"
$ws.Cells.Item(8, 5).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(8, 6).Value = "Shubman Gill "
$ws.Cells.Item(8, 7).Value = "70"
$ws.Cells.Item(8, 8).Value = "62"
$ws.Cells.Item(8, 9).Value = "5"
$ws.Cells.Item(8, 10).Value = "2"
$ws.Cells.Item(8, 11).Value = "112.90"

# Row 9
$ws.Cells.Item(9, 1).Value = " Sharjah"
$ws.Cells.Item(9, 2).Value = " October 03 2020"
$ws.Cells.Item(9, 3).Value = "Capitals won by 18 runs"
$ws.Cells.Item(9, 4).Value = "Kolkata Knight Riders"
$ws.Cells.Item(9, 5).Value = "Delhi Capitals"
$ws.Cells.Item(9, 6).Value = "Shubman Gill "
$ws.Cells.Item(9, 7).Value = "28"
$ws.Cells.Item(9, 8).Value = "22"
$ws.Cells.Item(9, 9).Value = "2"
$ws.Cells.Item(9, 10).Value = "1"
$ws.Cells.Item(9, 11).Value = "127.27"

# Row 10
$ws.Cells.Item(10, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(10, 2).Value = " October 29 2020"
$ws.Cells.Item(10, 3).Value = "Super Kings won by 6 wickets"
$ws.Cells.Item(10, 4).Value = "Kolkata Knight Riders"
$ws.Cells.Item(10, 5).Value = "Chennai Super Kings"
$ws.Cells.Item(10, 6).Value = "Shubman Gill "
$ws.Cells.Item(10, 7).Value = "26"
$ws.Cells.Item(10, 8).Value = "17"
$ws.Cells.Item(10, 9).Value = "4"
$ws.Cells.Item(10, 10).Value = "0"
$ws.Cells.Item(10, 11).Value = "152.94"

# Row 11
$ws.Cells.Item(11, 1).Value = " Sharjah"
$ws.Cells.Item(11, 2).Value = " October 26 2020"
$ws.Cells.Item(11, 3).Value = "Kings XI won by 8 wickets (with 7 balls remaining)"
$ws.Cells.Item(11, 4).Value = "Kolkata Knight Riders"
$ws.Cells.Item(11, 5).Value = "Kings XI Punjab"
$ws.Cells.Item(11, 6).Value = "Shubman Gill "
$ws.Cells.Item(11, 7).Value = "57"
$ws.Cells.Item(11, 8).Value = "45"
$ws.Cells.Item(11, 9).Value = "3"
$ws.Cells.Item(11, 10).Value = "4"
$ws.Cells.Item(11, 11).Value = "126.66"

# Row 12
$ws.Cells.Item(12, 1).Value = " Abu Dhabi"
$ws.Cells.Item(12, 2).Value = " October 21 2020"
$ws.Cells.Item(12, 3).Value = "RCB won by 8 wickets (with 39 balls remaining)"
$ws.Cells.Item(12, 4).Value = "Kolkata Knight Riders"
$ws.Cells.Item(12, 5).Value = "Royal Challengers Bangalore"
$ws.Cells.Item(12, 6).Value = "Shubman Gill "
$ws.Cells.Item(12, 7).Value = "1"
$ws.Cells.Item(12, 8).Value = "6"
$ws.Cells.Item(12, 9).Value = "0"
$ws.Cells.Item(12, 10).Value = "0"
$ws.Cells.Item(12, 11).Value = "16.66"

# Row 13
$ws.Cells.Item(13, 1).Value = " Abu Dhabi"
$ws.Cells.Item(13, 2).Value = " October 24 2020"
$ws.Cells.Item(13, 3).Value = "KKR won by 59 runs"
$ws.Cells.Item(13, 4).Value = "Kolkata Knight Riders"
$ws.Cells.Item(13, 5).Value = "Delhi Capitals"
$ws.Cells.Item(13, 6).Value = "Shubman Gill "
$ws.Cells.Item(13, 7).Value = "9"
$ws.Cells.Item(13, 8).Value = "8"
$ws.Cells.Item(13, 9).Value = "2"
$ws.Cells.Item(13, 10).Value = "0"
$ws.Cells.Item(13, 11).Value = "112.50"

# Row 14
$ws.Cells.Item(14, 1).Value = " Abu Dhabi"
$ws.Cells.Item(14, 2).Value = " September 23 2020"
$ws.Cells.Item(14, 3).Value = "Mumbai won by 49 runs"
$ws.Cells.Item(14, 4).Value = "Kolkata Knight Riders"
$ws.Cells.Item(14, 5).Value = "Mumbai Indians"
$ws.Cells.Item(14, 6).Value = "Shubman Gill "
$ws.Cells.Item(14, 7).Value = "7"
$ws.Cells.Item(14, 8).Value = "11"
$ws.Cells.Item(14, 9).Value = "1"
$ws.Cells.Item(14, 10).Value = "0"
$ws.Cells.Item(14, 11).Value = "63.63"

# Row 15
$ws.Cells.Item(15, 1).Value = " Sharjah"
$ws.Cells.Item(15, 2).Value = " October 12 2020"
$ws.Cells.Item(15, 3).Value = "RCB won by 82 runs"
$ws.Cells.Item(15, 4).Value = "Kolkata Knight Riders"
$ws.Cells.Item(15, 5).Value = "Royal Challengers Bangalore"
$ws.Cells.Item(15, 6).Value = "Shubman Gill "
$ws.Cells.Item(15, 7).Value = "34"
$ws.Cells.Item(15, 8).Value = "25"
$ws.Cells.Item(15, 9).Value = "3"
$ws.Cells.Item(15, 10).Value = "1"
$ws.Cells.Item(15, 11).Value = "136.00"
